# "ajout frais d'importation dans budget"
# Add a new budget line ("Importation batterie") on row 11 of the
# "Feuil1" budget sheet, and extend the MS sub-total formula (M7) to
# pick up the new row's reimbursable share (H11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 11 was a blank placeholder line; fill it in with the new
# "Importation batterie" budget entry.
#   C11 = Article, D11 = Quantité, E11 = Prix unitaire,
#   G11 = Achat par (MS), F11/H11 keep their pre-existing shared
#   formulas (Prix total = D*E, Valeur sur le robot = F/7) and will
#   recompute automatically once the inputs are set.
$ws.Range("C11").Value = "Importation batterie"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 34.8
$ws.Range("G11").Value = "MS"

# M7 ("MS :" sub-total) previously summed H6+H7; now also include the
# new row's reimbursable share (H11).
$ws.Range("M7").Formula = "=H6+H7+H11"

# Restore the author's final cell selection on this sheet.
[void]$ws.Range("M8").Select()
